# Update countries & provincias Spain
#
# The source data table on sheet "Pais" is kept sorted by total cases
# (column B) descending. A handful of countries had their case counts
# refreshed, which moves them past their neighbours in the ranking, and
# the "last updated" footer timestamp advances from 12:22 to 12:52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer: "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 12:52"

# Bielorrusia overtakes Japon (rows 32-33)
$ws.Cells.Item(32, 1).Value = "Bielorrusia"
$ws.Cells.Item(32, 2).Value = 14027
$ws.Cells.Item(32, 3).Value = 846
$ws.Cells.Item(32, 4).Value = 2386
$ws.Cells.Item(32, 5).Value = 11552
$ws.Cells.Item(32, 6).Value = 92
$ws.Cells.Item(32, 7).Value = 5
$ws.Cells.Item(32, 8).Value = 89

$ws.Cells.Item(33, 1).Value = "Japon"
$ws.Cells.Item(33, 2).Value = 13965
$ws.Cells.Item(33, 3).Value = 70
$ws.Cells.Item(33, 4).Value = 2368
$ws.Cells.Item(33, 5).Value = 11172
$ws.Cells.Item(33, 6).Value = 306
$ws.Cells.Item(33, 7).Value = 12
$ws.Cells.Item(33, 8).Value = 425

# Eslovenia refreshed in place (row 83)
$ws.Cells.Item(83, 2).Value = 1429
$ws.Cells.Item(83, 3).Value = 11
$ws.Cells.Item(83, 4).Value = 233
$ws.Cells.Item(83, 5).Value = 1105
$ws.Cells.Item(83, 7).Value = 2
$ws.Cells.Item(83, 8).Value = 91

# Albania overtakes Honduras (rows 95-96)
$ws.Cells.Item(95, 1).Value = "Albania"
$ws.Cells.Item(95, 2).Value = 773
$ws.Cells.Item(95, 3).Value = 7
$ws.Cells.Item(95, 4).Value = 470
$ws.Cells.Item(95, 5).Value = 272
$ws.Cells.Item(95, 6).Value = 4
$ws.Cells.Item(95, 7).Value = 1
$ws.Cells.Item(95, 8).Value = 31

$ws.Cells.Item(96, 1).Value = "Honduras"
$ws.Cells.Item(96, 2).Value = 771
$ws.Cells.Item(96, 3).Value = 33
$ws.Cells.Item(96, 4).Value = 79
$ws.Cells.Item(96, 5).Value = 621
$ws.Cells.Item(96, 6).Value = 10
$ws.Cells.Item(96, 7).Value = 5
$ws.Cells.Item(96, 8).Value = 71

# Vietnam refreshed in place (row 130)
$ws.Cells.Item(130, 4).Value = 219
$ws.Cells.Item(130, 5).Value = 51

# Sierra Leona overtakes Trinidad yTobago, Cabo Verde, Bermudas and Togo
# (rows 146-150 shift down one position each)
$ws.Cells.Item(146, 1).Value = "Sierra Leona"
$ws.Cells.Item(146, 2).Value = 116
$ws.Cells.Item(146, 3).Value = 12
$ws.Cells.Item(146, 4).Value = 14
$ws.Cells.Item(146, 5).Value = 98
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 4

$ws.Cells.Item(147, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(147, 2).Value = 116
$ws.Cells.Item(147, 3).Value = 0
$ws.Cells.Item(147, 4).Value = 72
$ws.Cells.Item(147, 5).Value = 36
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 0
$ws.Cells.Item(147, 8).Value = 8

$ws.Cells.Item(148, 1).Value = "Cabo Verde"
$ws.Cells.Item(148, 2).Value = 113
$ws.Cells.Item(148, 3).Value = 0
$ws.Cells.Item(148, 4).Value = 2
$ws.Cells.Item(148, 5).Value = 110
$ws.Cells.Item(148, 6).Value = 0
$ws.Cells.Item(148, 7).Value = 0
$ws.Cells.Item(148, 8).Value = 1

$ws.Cells.Item(149, 1).Value = "Bermudas"
$ws.Cells.Item(149, 2).Value = 111
$ws.Cells.Item(149, 3).Value = 0
$ws.Cells.Item(149, 4).Value = 48
$ws.Cells.Item(149, 5).Value = 57
$ws.Cells.Item(149, 6).Value = 10
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 6

$ws.Cells.Item(150, 1).Value = "Togo"
$ws.Cells.Item(150, 2).Value = 109
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(150, 4).Value = 64
$ws.Cells.Item(150, 5).Value = 38
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 7
